$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1
$ws.Range("E1").Value = "Phan Quan"
$ws.Range("B2").Value = "Trần Văn Đức"
$ws.Range("B3").Value = "20/6/2023 12:00:00 AM"
$ws.Range("B4").Value = "1000000.00 đồng"
$ws.Range("E4").Value = "50000.00 đồng"
$ws.Range("C5").Value = "Giá bán"
$ws.Range("A6").Value = "Phông trơn"
$ws.Range("B6").Value = "SWE"
$ws.Range("C6").Value = 150000
